$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Append facility rows that translate the dvdmt month/facility codes
# into their corresponding dhis2_name values (zone / district /
# dhis2_name / dvdmt_name columns A-D), continuing directly after
# the existing data which ends at row 764.
# -----------------------------------------------------------------

$data = @(
  @("South East", "Machinga", "Machinga HC",          "Machinga Health Centre"),
  @("South East", "Mangochi", "Monkey Bay",            "Monkey-Bay Community Hospital"),
  @("South East", "Mangochi", "St Martins",             "St Martins Molere Health Centre"),
  @("South East", "Mangochi", "Mtimabi",                "Mtimabii Health Centre"),
  @("South East", "Mangochi", "Mangochi Hosp",          "Mangochi District Hospital"),
  @("South East", "Mangochi", "Mulibwanji",             "Mulibwanji Hospital"),
  @("South East", "Mangochi", "katema",                 "Katema Health Centre"),
  @("South East", "Phalombe", "CHIRINGA CHAM",          "Chiringa Maternity"),
  @("South East", "Zomba",    "Namikango Maternity",    "Namikango Health Centre"),
  @("South East", "Zomba",    "Lungazi",                "Lungadzi Outreach Clinic"),
  @("South East", "Zomba",    "Elaine Zakresk PVT",     "Elaine Zakresh Outreach Clinic")
)

$startRow = 765
$endRow = $startRow + $data.Count - 1   # 775

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
}

# -----------------------------------------------------------------
# Row heights: every new row uses the 15.75pt custom row height.
# -----------------------------------------------------------------
$ws.Range("A$($startRow):D$endRow").RowHeight = 15.75

# -----------------------------------------------------------------
# Fonts / formatting for the new rows.
# -----------------------------------------------------------------

# Column A: Arial 10 (matches the workbook default look)
$rA = $ws.Range("A$($startRow):A$endRow")
$rA.Font.Name = "Arial"
$rA.Font.Size = 10

# Columns B & C: Calibri 12
$rBC = $ws.Range("B$($startRow):C$endRow")
$rBC.Font.Name = "Calibri"
$rBC.Font.Size = 12

# Column D, row 774: Calibri 11 (routed through Calibri 12 first so the
# 10pt "Calibri" intermediate state is shared with the B/C columns above)
$rD774 = $ws.Range("D774")
$rD774.Font.Name = "Calibri"
$rD774.Font.Size = 12
$rD774.Font.Size = 11

# Column D, remaining rows (765-773, 775): Courier New 14, green
$rD1 = $ws.Range("D765:D773")
$rD1.Font.Name = "Courier New"
$rD1.Font.Size = 14
$rD1.Font.Color = 32768
$rD1.Font.Family = 1

$rD2 = $ws.Range("D775")
$rD2.Font.Name = "Courier New"
$rD2.Font.Size = 14
$rD2.Font.Color = 32768
$rD2.Font.Family = 1

# -----------------------------------------------------------------
# Leave the view positioned near the newly appended rows.
# -----------------------------------------------------------------
$ws.Range("A735").Select() | Out-Null
$ws.Range("J777").Select() | Out-Null

Write-Host "Appended dhis2_name rows $startRow-$endRow"
